$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.789.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.635.98"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.860.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0768"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.83"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.794.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.31"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.78"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.58"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0493"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.25"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.129.80"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.548"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.60"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.37"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.808"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.769.27"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.24"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0506"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.417"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.53"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.15%  "
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.17%  "
